# Adds the "Manipulacion y creacion de columnas" content to the end of the
# document body (right before the final section break), as a single batch of
# 14 new paragraphs, by inserting a WordprocessingML package fragment at the
# end of the document's main story range.
$d = $word.ActiveDocument
$r = $d.Content
$r.Collapse(0) | Out-Null

$xmlFragment = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Para reordenar las columnas de un dataframe se hace de esta forma:duplicate[['ID','Nombre','Edad', 'Salario','Es duplicado?','Fecha_Ingreso']]</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    ID  Nombre  Edad  Salario Es duplicado? Fecha_Ingreso</w:t><w:br/><w:t>0    1     Ana    23    50000            No    2020-01-10</w:t><w:br/><w:t>1    2    Luis    35    60000            No    2019-06-15</w:t><w:br/><w:t>2    3  Carlos    30    58000            No    2020-03-25</w:t><w:br/><w:t>3    4   María    22    58000            No    2018-07-30</w:t><w:br/><w:t>4    5   Elena    28    49000            No    2018-07-30</w:t><w:br/><w:t>5    3  Carlos    30    58000            No    2020-03-25</w:t><w:br/><w:t>6    7   Sofía    34    72000            No    2016-11-22</w:t><w:br/><w:t>7    8  Andrés    45    61000            No    2021-04-01</w:t><w:br/><w:t>8    9   Laura    29    52000            No    2015-05-19</w:t><w:br/><w:t>9   10     Ana    41    67000            No    2014-08-25</w:t><w:br/><w:t>10   2    Luis    35    60000            Si    2019-06-15</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Para elegir columnas especificas de un dataframe se puede usa la propiedad .loc de esta manera: duplicate.loc[:,['ID','Nombre']] </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    ID  Nombre</w:t><w:br/><w:t>0    1     Ana</w:t><w:br/><w:t>1    2    Luis</w:t><w:br/><w:t>2    3  Carlos</w:t><w:br/><w:t>3    4   María</w:t><w:br/><w:t>4    5   Elena</w:t><w:br/><w:t>5    3  Carlos</w:t><w:br/><w:t>6    7   Sofía</w:t><w:br/><w:t>7    8  Andrés</w:t><w:br/><w:t>8    9   Laura</w:t><w:br/><w:t>9   10     Ana</w:t><w:br/><w:t>10   2    Luis</w:t></w:r></w:p><w:p><w:r><w:t>El primer parametro de la propiedad .loc se utiliza para especificar el rango de las filas que se selecionaran, ejemplo de filas del 1 al 9 / 1:9</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">   ID  Nombre</w:t><w:br/><w:t>1   2    Luis</w:t><w:br/><w:t>2   3  Carlos</w:t><w:br/><w:t>3   4   María</w:t><w:br/><w:t>4   5   Elena</w:t><w:br/><w:t>5   3  Carlos</w:t><w:br/><w:t>6   7   Sofía</w:t><w:br/><w:t>7   8  Andrés</w:t><w:br/><w:t>8   9   Laura</w:t></w:r></w:p><w:p><w:r><w:t>Para eliminar una columna, ejemplo la duplicado, se utiliza el metodo ..drop(columns=['Es duplicado?'])</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    ID  Nombre  Edad  Salario Fecha_Ingreso</w:t><w:br/><w:t>0    1     Ana    23    50000    2020-01-10</w:t><w:br/><w:t>1    2    Luis    35    60000    2019-06-15</w:t><w:br/><w:t>2    3  Carlos    30    58000    2020-03-25</w:t><w:br/><w:t>3    4   María    22    58000    2018-07-30</w:t><w:br/><w:t>4    5   Elena    28    49000    2018-07-30</w:t><w:br/><w:t>5    3  Carlos    30    58000    2020-03-25</w:t><w:br/><w:t>6    7   Sofía    34    72000    2016-11-22</w:t><w:br/><w:t>7    8  Andrés    45    61000    2021-04-01</w:t><w:br/><w:t>8    9   Laura    29    52000    2015-05-19</w:t><w:br/><w:t>9   10     Ana    41    67000    2014-08-25</w:t><w:br/><w:t>10   2    Luis    35    60000    2019-06-15</w:t></w:r></w:p><w:p><w:r><w:t>Si quiero filtra los salarios superiores a 50,000 se hace de la siguiente manera: duplicate.loc[duplicate['Salario' ]&gt; 50000]</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    ID  Nombre  Edad  Salario Fecha_Ingreso</w:t><w:br/><w:t>1    2    Luis    35    60000    2019-06-15</w:t><w:br/><w:t>2    3  Carlos    30    58000    2020-03-25</w:t><w:br/><w:t>3    4   María    22    58000    2018-07-30</w:t><w:br/><w:t>5    3  Carlos    30    58000    2020-03-25</w:t><w:br/><w:t>6    7   Sofía    34    72000    2016-11-22</w:t><w:br/><w:t>7    8  Andrés    45    61000    2021-04-01</w:t><w:br/><w:t>8    9   Laura    29    52000    2015-05-19</w:t><w:br/><w:t>9   10     Ana    41    67000    2014-08-25</w:t><w:br/><w:t>10   2    Luis    35    60000    2019-06-15</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Para agregar una columna nueva a un dataframe se hace de esta manera:duplicate['Posición'] agregando los valores que tendra la columna: </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    ID  Nombre  Edad  Salario Fecha_Ingreso Posición</w:t><w:br/><w:t>0    1     Ana    23    50000    2020-01-10   junior</w:t><w:br/><w:t>1    2    Luis    35    60000    2019-06-15   junior</w:t><w:br/><w:t>2    3  Carlos    30    58000    2020-03-25   junior</w:t><w:br/><w:t>3    4   María    22    58000    2018-07-30   junior</w:t><w:br/><w:t>4    5   Elena    28    49000    2018-07-30   junior</w:t><w:br/><w:t>5    3  Carlos    30    58000    2020-03-25   junior</w:t><w:br/><w:t>6    7   Sofía    34    72000    2016-11-22   Senior</w:t><w:br/><w:t>7    8  Andrés    45    61000    2021-04-01      mid</w:t><w:br/><w:t>8    9   Laura    29    52000    2015-05-19   junior</w:t><w:br/><w:t>9   10     Ana    41    67000    2014-08-25   Senior</w:t><w:br/><w:t>10   2    Luis    35    60000    2019-06-15   junior</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Esta nueva columna se calcula cuanto se le descuenta de afp + ars: </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    ID  Nombre  Edad  Salario Fecha_Ingreso Posición  AFP + ARS  Salario_Neto</w:t><w:br/><w:t>0    1     Ana    23    50000    2020-01-10   junior     2955.0       47045.0</w:t><w:br/><w:t>1    2    Luis    35    60000    2019-06-15   junior     3546.0       56454.0</w:t><w:br/><w:t>2    3  Carlos    30    58000    2020-03-25   junior     3427.8       54572.2</w:t><w:br/><w:t>3    4   María    22    58000    2018-07-30   junior     3427.8       54572.2</w:t><w:br/><w:t>4    5   Elena    28    49000    2018-07-30   junior     2895.9       46104.1</w:t><w:br/><w:t>5    3  Carlos    30    58000    2020-03-25   junior     3427.8       54572.2</w:t><w:br/><w:t>6    7   Sofía    34    72000    2016-11-22   Senior     4255.2       67744.8</w:t><w:br/><w:t>7    8  Andrés    45    61000    2021-04-01      mid     3605.1       57394.9</w:t><w:br/><w:t>8    9   Laura    29    52000    2015-05-19   junior     3073.2       48926.8</w:t><w:br/><w:t>9   10     Ana    41    67000    2014-08-25   Senior     3959.7       63040.3</w:t><w:br/><w:t>10   2    Luis    35    60000    2019-06-15   junior     3546.0       56454.0</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$r.InsertXML($xmlFragment)
